$wb = $excel.ActiveWorkbook

# --- Sheet3 tab (physical sheet2.xml): add two new values, keep existing D10 ---
$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet3.Range("G1").Value = 31
$sheet3.Range("A4").Value = 56

# --- Sheet4 tab (physical sheet4.xml -> sheet5.xml): add three new values, keep existing E4 ---
$sheet4 = $wb.Worksheets.Item("Sheet4")
$sheet4.Range("B3").Value = 52
$sheet4.Range("G6").Value = 29
$sheet4.Range("D7").Value = 22

# --- Insert a brand new sheet ("Sheet5") right after "Sheet3"; becomes the active tab ---
$sheet3 = $wb.Worksheets.Item("Sheet3")
$newSheet = $wb.Worksheets.Add($null, $sheet3)
$newSheet.Range("A1").Value = 31
$newSheet.Range("B1").Value = 53
$newSheet.Range("C1").Value = 31

# --- Re-fetch worksheet references (inserting a sheet can shift cached index-bound
#     references) and set final selections to match the target view state ---
$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet3.Range("I11").Select()

$sheet4 = $wb.Worksheets.Item("Sheet4")
$sheet4.Range("G6").Select()

$newSheet = $wb.Worksheets.Item("Sheet5")
$newSheet.Range("C1").Select()
